# Regenerate the "K" column (column G) values for the save_data sheet.
# Commit: "regen save_data to use K instead of Strike#, regen std/mean,
#          calc and write s_vals"
#
# Only column G (header "K" in row 1) changes; every other column is left
# untouched. Rows 6 and 14 already hold the correct (0) value so they are
# not listed among the updates below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 1
    4  = 0
    5  = 0
    7  = 0
    8  = 1
    9  = 0
    10 = 1
    11 = 3
    12 = 0
    13 = 1
    15 = 2
    16 = 1
    17 = 0
    18 = 0
    19 = 1
    20 = 1
    21 = 2
    22 = 2
    23 = 0
    24 = 0
    25 = 1
    26 = 2
    27 = 0
    28 = 2
    29 = 1
    30 = 1
    31 = 2
    32 = 3
    33 = 1
    34 = 2
    35 = 4
    36 = 1
    37 = 1
    38 = 2
    39 = 1
    40 = 1
    41 = 1
    42 = 0
    43 = 0
    44 = 1
    45 = 3
    46 = 1
    47 = 1
    48 = 1
    49 = 1
    50 = 2
    51 = 2
    52 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
